$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.282.44"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'1.929.98"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'248.87"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'0.7162"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.3211"
$ws.Range("E8").Value = "  -4.37%  "
$ws.Range("D9").Value = "'27.71"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("D10").Value = "'0.07101"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").Value = "'0.7917"
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "'0.08003"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'1.933.61"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'5.397"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "'94.88"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'14.66"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "'30.300.05"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'256.92"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "'0.000008051"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "'5.769"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").Value = "'2.184.70"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.836"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "'9.548"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'164.82"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").Value = "'19.12"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'2.274"
$ws.Range("E28").Value = "  -6.07%  "
$ws.Range("D29").Value = "'0.1271"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("D30").Value = "'1.355"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "'4.398"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "'4.132"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("D34").Value = "'0.05146"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "'1.269"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "'0.7449"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'2.764"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "'0.01963"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "'2.798"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "'78.75"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'6.361"
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").Value = "'0.4514"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'1.997"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'100.50"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'9.772"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'7.440"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'36.70"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "'946.23"
$ws.Range("E50").Value = "  +9.37%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06108"
$ws.Range("E51").Value = "  +1.77%  "
